$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 — copy formatting from the existing
# header cell (G1: bold, centered, bordered) so H1 matches the rest of
# row 1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1").Value = "Save"

# New "Save" data value for the single data row.
$ws.Range("H2").Value = 0
